$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 11: Y / Test Testov / Rupee
$ws.Range("A11").Value = "Y"
$ws.Range("B11").Value = "Test Testov"
$ws.Range("C11").Value = "Rupee"

# Row 12: N / Kalim Neon / Dollar
$ws.Range("A12").Value = "N"
$ws.Range("B12").Value = "Kalim Neon"
$ws.Range("C12").Value = "Dollar"

# Move the active selection the way Excel would after entering the last row
$ws.Range("C13").Select()
